$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.590.56'
$ws.Range("E2").Value = '  -2.96%  '

$ws.Range("D3").Value = '2.625.21'
$ws.Range("E3").Value = '  -1.08%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''573.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.58%  '

$ws.Range("D6").Value = '''154.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.96%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").Value = '''0.621'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.00%  '

$ws.Range("D9").Value = '2.622.18'
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").Value = '''0.116'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.51%  '

$ws.Range("D11").Value = '''5.76'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.69%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.379'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.12%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.156'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.14%  '

$ws.Range("D14").Value = '''28.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").Value = '3.092.10'
$ws.Range("E15").Value = '  -1.09%  '

$ws.Range("D16").Value = '''0.0000184'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.92%  '

$ws.Range("D17").Value = '63.495.70'
$ws.Range("E17").Value = '  -2.99%  '

$ws.Range("D18").Value = '2.629.81'
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").Value = '''12.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.47%  '

$ws.Range("D20").Value = '''7.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").Value = '''4.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.69%  '

$ws.Range("D22").Value = '''341.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.45%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").Value = '''67.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.43%  '

$ws.Range("D25").Value = '''1.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.10%  '

$ws.Range("D26").Value = '''0.0000107'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.92%  '

$ws.Range("D27").Value = '''583.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.08%  '

$ws.Range("D28").Value = '''9.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.71%  '

$ws.Range("D29").Value = '''1.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.50%  '

$ws.Range("D30").Value = '''0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("D31").Value = '''0.160'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.32%  '

$ws.Range("D32").Value = '''7.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.41%  '

$ws.Range("D33").Value = '''2.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("D34").Value = '''1.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.57%  '

$ws.Range("D35").Value = '''6.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").Value = '''5.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.23%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '''0.399'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.83%  '

$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = '''0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").Value = '''19.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.72%  '

$ws.Range("D40").Value = '''153.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.65%  '

$ws.Range("D41").Value = '''1.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.86%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").Value = '''41.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("D44").Value = '''2.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.91%  '

$ws.Range("D45").Value = '''158.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''23.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.31%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '''3.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.36%  '

$ws.Range("D48").Value = '''0.0582'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.50%  '

$ws.Range("D49").Value = '''0.629'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.16%  '

$ws.Range("D50").Value = '''0.0993'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").Value = '''0.0246'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.13%  '

